$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 63.91118233333333
$ws.Range("H2").Value = 191.733547
$ws.Range("I2").Value = 0.4067926910433548
$ws.Range("J2").Value = 0.4067926910433549
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.051093
$ws.Range("N2").Value = 0.153279
$ws.Range("O2").Value = 0.01450579975525089
$ws.Range("P2").Value = 0.01450579975525089
$ws.Range("Q2").Value = 3.265414038957
$ws.Range("R2").Value = 29.388726350613
$ws.Range("S2").Value = 0.005900853318174545
$ws.Range("T2").Value = 0.005900853318174548

$ws.Range("G3").Value = 63.91118233333333
$ws.Range("H3").Value = 191.733547
$ws.Range("I3").Value = 0.4067926910433548
$ws.Range("J3").Value = 0.4067926910433549
$ws.Range("O3").Value = 0.2313022967634575
$ws.Range("P3").Value = 0.2313022967634575
$ws.Range("Q3").Value = 52.06867458796854
$ws.Range("R3").Value = 468.6180712917169
$ws.Range("S3").Value = 0.09409208374491551
$ws.Range("T3").Value = 0.09409208374491555

$ws.Range("G4").Value = 63.91118233333333
$ws.Range("H4").Value = 191.733547
$ws.Range("I4").Value = 0.4067926910433548
$ws.Range("J4").Value = 0.4067926910433549
$ws.Range("M4").Value = 2.656449666666667
$ws.Range("N4").Value = 7.969348999999999
$ws.Range("O4").Value = 0.7541919034812916
$ws.Range("P4").Value = 0.7541919034812917
$ws.Range("Q4").Value = 169.7768390056559
$ws.Range("R4").Value = 1527.991551050903
$ws.Range("S4").Value = 0.3067997539802647
$ws.Range("T4").Value = 0.3067997539802648

$ws.Range("I5").Value = 0.3656254573230189
$ws.Range("J5").Value = 0.365625457323019
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.051093
$ws.Range("N5").Value = 0.153279
$ws.Range("O5").Value = 0.01450579975525089
$ws.Range("P5").Value = 0.01450579975525089
$ws.Range("Q5").Value = 2.9349556362
$ws.Range("R5").Value = 26.4146007258
$ws.Range("S5").Value = 0.005303689669349741
$ws.Range("T5").Value = 0.005303689669349743

$ws.Range("I6").Value = 0.3656254573230189
$ws.Range("J6").Value = 0.365625457323019
$ws.Range("O6").Value = 0.2313022967634575
$ws.Range("P6").Value = 0.2313022967634575
$ws.Range("Q6").Value = 46.79934860579999
$ws.Range("S6").Value = 0.08457000803400377
$ws.Range("T6").Value = 0.0845700080340038

$ws.Range("I7").Value = 0.3656254573230189
$ws.Range("J7").Value = 0.365625457323019
$ws.Range("M7").Value = 2.656449666666667
$ws.Range("N7").Value = 7.969348999999999
$ws.Range("O7").Value = 0.7541919034812916
$ws.Range("P7").Value = 0.7541919034812917
$ws.Range("Q7").Value = 152.5955007822
$ws.Range("R7").Value = 1373.3595070398
$ws.Range("S7").Value = 0.2757517596196654
$ws.Range("T7").Value = 0.2757517596196655

$ws.Range("G8").Value = 35.755375
$ws.Range("H8").Value = 107.266125
$ws.Range("I8").Value = 0.2275818516336261
$ws.Range("J8").Value = 0.2275818516336262
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.051093
$ws.Range("N8").Value = 0.153279
$ws.Range("O8").Value = 0.01450579975525089
$ws.Range("P8").Value = 0.01450579975525089
$ws.Range("Q8").Value = 1.826849374875
$ws.Range("R8").Value = 16.441644373875
$ws.Range("S8").Value = 0.003301256767726597
$ws.Range("T8").Value = 0.003301256767726599

$ws.Range("G9").Value = 35.755375
$ws.Range("H9").Value = 107.266125
$ws.Range("I9").Value = 0.2275818516336261
$ws.Range("J9").Value = 0.2275818516336262
$ws.Range("O9").Value = 0.2313022967634575
$ws.Range("P9").Value = 0.2313022967634575
$ws.Range("Q9").Value = 29.13003511554167
$ws.Range("R9").Value = 262.170316039875
$ws.Range("S9").Value = 0.05264020498453813
$ws.Range("T9").Value = 0.05264020498453816

$ws.Range("G10").Value = 35.755375
$ws.Range("H10").Value = 107.266125
$ws.Range("I10").Value = 0.2275818516336261
$ws.Range("J10").Value = 0.2275818516336262
$ws.Range("M10").Value = 2.656449666666667
$ws.Range("N10").Value = 7.969348999999999
$ws.Range("O10").Value = 0.7541919034812916
$ws.Range("P10").Value = 0.7541919034812917
$ws.Range("Q10").Value = 94.98235400029166
$ws.Range("R10").Value = 854.8411860026249
$ws.Range("S10").Value = 0.1716403898813614
$ws.Range("T10").Value = 0.1716403898813614
